# Update countries & provincias Spain
#
# The underlying dataset is sorted by total cases (column B) descending.
# A handful of countries' case counts were refreshed, which changed their
# rank/position in the table, so several rows swap / shift places and some
# rows simply get refreshed figures. The timestamp banner in A1 is bumped
# too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($row, $country, $b, $c, $d, $e, $f, $g, $h)
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Peru overtakes Mexico (rows 9/10 swap, Peru gets refreshed totals)
Set-Row 9  "Peru"   516296 0    354232 136208 0 0   25856
Set-Row 10 "Mexico" 511369 5618 345653 109808 0 615 55908

# Australia: active/recovered refreshed slightly (no rank change)
$ws.Cells.Item(71, 4).Value = 13355
$ws.Cells.Item(71, 5).Value = 9301

# Paraguay overtakes Zambia and Guayana Francesa (rows 89/90/91 shift)
Set-Row 89 "Paraguay"         9022 0 5657 3257 0 0 108
Set-Row 90 "Zambia"           9021 0 7586 1179 0 0 256
Set-Row 91 "Guayana Francesa" 8549 0 7841  655 0 0  53

# Belice: minor figures refreshed (no rank change)
$ws.Cells.Item(174, 5).Value = 321
$ws.Cells.Item(174, 7).Value = 1
$ws.Cells.Item(174, 8).Value = 3

# Islas Turcas y Caicos overtakes Camboya, Papua Nueva Guinea and
# San Martin (Parte Holandesa) (rows 180/181/182/183 shift)
Set-Row 180 "Islas Turcas y Caicos"         274 16 54  218 0 0  2
Set-Row 181 "Camboya"                       273 0  225  48 0 0  0
Set-Row 182 "Papua Nueva Guinea"            271 0   78 190 0 0  3
Set-Row 183 "San Martin (Parte Holandesa)"  269 6  102 150 0 0 17

# Montserrat overtakes Islas Malvinas (rows 213/214 swap)
Set-Row 213 "Montserrat"      13 0 12 0 0 0 1
Set-Row 214 "Islas Malvinas"  13 0 13 0 0 0 0

# Bump the "last updated" banner
$ws.Range("A1").Value = "Datos actualizados a 15 de Agosto de 2020 a las 05:30"
